# Feature : 화살 ObjectPooling 구현
# Rework the Range weapon's projectile reference from a literal resource
# path (string used to load/instantiate a prefab each shot) to a pooling
# tag (string key used to fetch an already-pooled instance).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Column header: "_projectilePath" -> "_projectileTag"
$ws.Range("Q1").Value = "_projectileTag"

# Bow's (row 6) projectile reference now stores the pool tag "Arrow"
# instead of the old "Resources/Items/Prefabs/Projectiles/" prefab path.
$ws.Range("Q6").Value = "Arrow"

# Leave the cursor where the edit happened, like the authored workbook.
$ws.Range("Q7").Select()
